$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay text (matches the
# pre-existing inline-string typing of these columns); force Text format
# before assignment so Excel does not auto-coerce them to numbers.
$textCells = @('D4', 'D5', 'D6', 'D10', 'D11', 'D14', 'D18', 'D20', 'D22', 'D23', 'D24', 'D25', 'D28', 'D29', 'D30', 'D31', 'D32', 'D35', 'D37', 'D38', 'D39', 'D40', 'D42', 'D44', 'D45', 'D46', 'D47', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.617.93'
$ws.Range('E2').Value = '  +3.48%  '
$ws.Range('D3').Value = '2.257.67'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '303.67'
$ws.Range('E5').Value = '  +2.41%  '
$ws.Range('D6').Value = '91.55'
$ws.Range('E6').Value = '  +4.36%  '
$ws.Range('E7').Value = '  +2.61%  '
$ws.Range('E9').Value = '  +1.53%  '
$ws.Range('D10').Value = '32.08'
$ws.Range('E10').Value = '  +3.84%  '
$ws.Range('D11').Value = '52.88'
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('E12').Value = '  +1.88%  '
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').Value = '6.56'
$ws.Range('E14').Value = '  +2.47%  '
$ws.Range('D15').Value = '2.606.87'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').Value = '2.277.15'
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('D18').Value = '0.760'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').Value = '41.528.75'
$ws.Range('E19').Value = '  +3.51%  '
$ws.Range('D20').Value = '12.37'
$ws.Range('E20').Value = '  +9.20%  '
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('D22').Value = '5.90'
$ws.Range('E22').Value = '  +2.22%  '
$ws.Range('D23').Value = '66.61'
$ws.Range('E23').Value = '  +1.34%  '
$ws.Range('D24').Value = '240.04'
$ws.Range('E24').Value = '  +1.83%  '
$ws.Range('D25').Value = '2.59'
$ws.Range('E25').Value = '  +3.77%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  +5.60%  '
$ws.Range('D28').Value = '23.92'
$ws.Range('E28').Value = '  +2.79%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.13'
$ws.Range('E29').Value = '  +2.65%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '9.49'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('D31').Value = '159.79'
$ws.Range('E31').Value = '  +2.23%  '
$ws.Range('D32').Value = '34.12'
$ws.Range('E32').Value = '  +6.14%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +3.94%  '
$ws.Range('D35').Value = '0.0741'
$ws.Range('E35').Value = '  +3.71%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '2.37'
$ws.Range('E37').Value = '  +1.97%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.104'
$ws.Range('E38').Value = '  +3.15%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '0.116'
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('D40').Value = '16.56'
$ws.Range('E40').Value = '  +6.36%  '
$ws.Range('E41').Value = '  +3.46%  '
$ws.Range('D42').Value = '3.93'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').Value = '2.046.44'
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').Value = '19.50'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '10.36'
$ws.Range('E45').Value = '  +3.70%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0277'
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('D47').Value = '2.85'
$ws.Range('E47').Value = '  +0.96%  '
$ws.Range('E48').Value = '  +7.00%  '
$ws.Range('E49').Value = '  +4.81%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '1.16'
$ws.Range('E50').Value = '  +2.44%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '72.59'
$ws.Range('E51').Value = '  +6.69%  '
